# Обновления апрельского релиза: добавление новых пар "Имена"/"Перевод"
# (Add files via upload) — append 7 new name translation rows to the
# "Имена"/"Перевод" table (columns M/N) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M18").Value = "Ailish"
$ws.Range("N18").Value = "Айлиш"

$ws.Range("M19").Value = "Hrafnkel"
$ws.Range("N19").Value = "Храфнкель"

$ws.Range("M20").Value = "Cuana"
$ws.Range("N20").Value = "Куана"

$ws.Range("M21").Value = "Dylan"
$ws.Range("N21").Value = "Дилан"

$ws.Range("M22").Value = "Atiq"
$ws.Range("N22").Value = "Атик"

$ws.Range("M23").Value = "Meaghan"
$ws.Range("N23").Value = "Миган"

$ws.Range("M24").Value = "Maebh"
$ws.Range("N24").Value = "Мэйб"
